$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly record at row 355 (shifts existing rows 355-401
# down to 356-402), matching the "Fruta / hortaliza, semanal" commit.
$ws.Rows("355").Insert()

$ws.Cells.Item(355, 1).Value = 1
$ws.Cells.Item(355, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(355, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(355, 4).Value = 45077
$ws.Cells.Item(355, 5).Value = 15
$ws.Cells.Item(355, 6).Value = "Fruta"
$ws.Cells.Item(355, 7).Value = 100108
$ws.Cells.Item(355, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(355, 9).Value = 100108006
$ws.Cells.Item(355, 10).Value = "Plátano"
$ws.Cells.Item(355, 11).Value = "Sin especificar"
$ws.Cells.Item(355, 12).Value = "Pintón"
$ws.Cells.Item(355, 13).Value = 200
$ws.Cells.Item(355, 14).Value = 19000
$ws.Cells.Item(355, 15).Value = 20000
$ws.Cells.Item(355, 16).Value = 19500
$ws.Cells.Item(355, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(355, 18).Value = "Ecuador"
$ws.Cells.Item(355, 19).Value = 975
$ws.Cells.Item(355, 20).Value = 20
